# Scheduled market-data refresh: updates currentAveragePrice* / Leve
# price / profit columns (H:N) across the job sheets, per the latest
# Universalis pull. Leve name/item/level/exp/gil/amount/item-id columns
# (A:G) are untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1662.1818
$ws.Cells.Item(28, 9).Value = 1288
$ws.Cells.Item(28, 10).Value = 2660
$ws.Cells.Item(28, 11).Value = 1288
$ws.Cells.Item(28, 12).Value = 2660
$ws.Cells.Item(28, 13).Value = -803
$ws.Cells.Item(28, 14).Value = -3630

$ws.Cells.Item(33, 8).Value = 132.22223
$ws.Cells.Item(33, 9).Value = 132.22223
$ws.Cells.Item(33, 11).Value = 132.22223
$ws.Cells.Item(33, 13).Value = 96.77777

$ws.Cells.Item(40, 8).Value = 10874.75
$ws.Cells.Item(40, 9).Value = 9600
$ws.Cells.Item(40, 10).Value = 12345.615
$ws.Cells.Item(40, 11).Value = 9600
$ws.Cells.Item(40, 12).Value = 12345.615
$ws.Cells.Item(40, 13).Value = -9425
$ws.Cells.Item(40, 14).Value = -12695.615

$ws.Cells.Item(64, 8).Value = 6214950
$ws.Cells.Item(64, 9).Value = 21741628
$ws.Cells.Item(64, 10).Value = 4279.4
$ws.Cells.Item(64, 11).Value = 21741628
$ws.Cells.Item(64, 12).Value = 4279.4
$ws.Cells.Item(64, 13).Value = -21741380
$ws.Cells.Item(64, 14).Value = -4775.4

$ws.Cells.Item(67, 8).Value = 6214950
$ws.Cells.Item(67, 9).Value = 21741628
$ws.Cells.Item(67, 10).Value = 4279.4
$ws.Cells.Item(67, 11).Value = 21741628
$ws.Cells.Item(67, 12).Value = 4279.4
$ws.Cells.Item(67, 13).Value = -21740770
$ws.Cells.Item(67, 14).Value = -5995.4

$ws.Cells.Item(80, 8).Value = 2076
$ws.Cells.Item(80, 9).Value = 1900
$ws.Cells.Item(80, 10).Value = 2105.3333
$ws.Cells.Item(80, 11).Value = 5700
$ws.Cells.Item(80, 12).Value = 6315.999899999999
$ws.Cells.Item(80, 13).Value = -4702
$ws.Cells.Item(80, 14).Value = -8311.999899999999

$ws.Cells.Item(83, 8).Value = 2076
$ws.Cells.Item(83, 9).Value = 1900
$ws.Cells.Item(83, 10).Value = 2105.3333
$ws.Cells.Item(83, 11).Value = 17100
$ws.Cells.Item(83, 12).Value = 18947.9997
$ws.Cells.Item(83, 13).Value = -12108
$ws.Cells.Item(83, 14).Value = -28931.9997

$ws.Cells.Item(86, 8).Value = 1159880.4
$ws.Cells.Item(86, 9).Value = 1941959.5
$ws.Cells.Item(86, 11).Value = 1941959.5
$ws.Cells.Item(86, 13).Value = -1940836.5

$ws.Cells.Item(88, 8).Value = 12510118
$ws.Cells.Item(88, 10).Value = 10914.077
$ws.Cells.Item(88, 12).Value = 10914.077
$ws.Cells.Item(88, 14).Value = -11726.077

$ws.Cells.Item(89, 8).Value = 1159880.4
$ws.Cells.Item(89, 9).Value = 1941959.5
$ws.Cells.Item(89, 11).Value = 9709797.5
$ws.Cells.Item(89, 13).Value = -9704181.5

$ws.Cells.Item(91, 8).Value = 12510118
$ws.Cells.Item(91, 10).Value = 10914.077
$ws.Cells.Item(91, 12).Value = 10914.077
$ws.Cells.Item(91, 14).Value = -13722.077

$ws.Cells.Item(137, 8).Value = 12266366
$ws.Cells.Item(137, 9).Value = 910895.6
$ws.Cells.Item(137, 11).Value = 2732686.8
$ws.Cells.Item(137, 13).Value = -2730136.8

$ws.Cells.Item(138, 8).Value = 3630.394
$ws.Cells.Item(138, 9).Value = 1342.9565
$ws.Cells.Item(138, 10).Value = 4322.6445
$ws.Cells.Item(138, 11).Value = 4028.8695
$ws.Cells.Item(138, 12).Value = 12967.9335
$ws.Cells.Item(138, 13).Value = 1111.1305
$ws.Cells.Item(138, 14).Value = -23247.9335

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 12958.098
$ws.Cells.Item(32, 9).Value = 12281.049
$ws.Cells.Item(32, 11).Value = 12281.049
$ws.Cells.Item(32, 13).Value = -11994.049

$ws.Cells.Item(110, 8).Value = 1278325.8
$ws.Cells.Item(110, 9).Value = 1702434.4
$ws.Cells.Item(110, 11).Value = 1702434.4
$ws.Cells.Item(110, 13).Value = -1700389.4

$ws.Cells.Item(132, 8).Value = 14954.367
$ws.Cells.Item(132, 9).Value = 23810.424
$ws.Cells.Item(132, 11).Value = 71431.272
$ws.Cells.Item(132, 13).Value = -68901.272

$ws.Cells.Item(135, 8).Value = 74809.664
$ws.Cells.Item(135, 10).Value = 74809.664
$ws.Cells.Item(135, 12).Value = 74809.664
$ws.Cells.Item(135, 14).Value = -84949.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1303504.2
$ws.Cells.Item(99, 9).Value = 1489433.4
$ws.Cells.Item(99, 11).Value = 1489433.4
$ws.Cells.Item(99, 13).Value = -1487935.4

$ws.Cells.Item(134, 8).Value = 1136
$ws.Cells.Item(134, 9).Value = 1119.5555
$ws.Cells.Item(134, 11).Value = 3358.6665
$ws.Cells.Item(134, 13).Value = -823.6664999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 242.94444
$ws.Cells.Item(7, 9).Value = 214
$ws.Cells.Item(7, 10).Value = 261.36365
$ws.Cells.Item(7, 11).Value = 214
$ws.Cells.Item(7, 12).Value = 261.36365
$ws.Cells.Item(7, 13).Value = -101
$ws.Cells.Item(7, 14).Value = -487.36365

$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 14).ClearContents()

$ws.Cells.Item(99, 8).Value = 7631.8423
$ws.Cells.Item(99, 10).Value = 7987.067
$ws.Cells.Item(99, 12).Value = 7987.067
$ws.Cells.Item(99, 14).Value = -10983.067

$ws.Cells.Item(107, 8).Value = 2020834.1
$ws.Cells.Item(107, 10).Value = 618
$ws.Cells.Item(107, 12).Value = 618
$ws.Cells.Item(107, 14).Value = -4458

$ws.Cells.Item(126, 8).Value = 7631.8423
$ws.Cells.Item(126, 10).Value = 7987.067
$ws.Cells.Item(126, 12).Value = 23961.201
$ws.Cells.Item(126, 14).Value = -28901.201

$ws.Cells.Item(141, 8).Value = 81670.55499999999
$ws.Cells.Item(141, 10).Value = 82927.88
$ws.Cells.Item(141, 12).Value = 82927.88
$ws.Cells.Item(141, 14).Value = -93287.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 3627184.8
$ws.Cells.Item(4, 9).Value = 895128.2
$ws.Cells.Item(4, 11).Value = 2685384.6
$ws.Cells.Item(4, 13).Value = -2685272.6

$ws.Cells.Item(14, 8).Value = 251.23077
$ws.Cells.Item(14, 9).Value = 251.23077
$ws.Cells.Item(14, 11).Value = 753.69231
$ws.Cells.Item(14, 13).Value = -580.69231

$ws.Cells.Item(131, 8).Value = 12486111
$ws.Cells.Item(131, 10).Value = 13694351
$ws.Cells.Item(131, 12).Value = 41083053
$ws.Cells.Item(131, 14).Value = -41093133

$ws.Cells.Item(136, 8).Value = 11205.782
$ws.Cells.Item(136, 9).Value = 5621.8
$ws.Cells.Item(136, 10).Value = 12756.889
$ws.Cells.Item(136, 11).Value = 16865.4
$ws.Cells.Item(136, 12).Value = 38270.667
$ws.Cells.Item(136, 13).Value = -11765.4
$ws.Cells.Item(136, 14).Value = -48470.667

$ws.Cells.Item(138, 8).Value = 5093.727
$ws.Cells.Item(138, 9).Value = 10555.5
$ws.Cells.Item(138, 10).Value = 3880
$ws.Cells.Item(138, 11).Value = 31666.5
$ws.Cells.Item(138, 12).Value = 11640
$ws.Cells.Item(138, 13).Value = -26526.5
$ws.Cells.Item(138, 14).Value = -21920

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 1296976.4
$ws.Cells.Item(80, 9).Value = 2379099.5
$ws.Cells.Item(80, 10).Value = 34499.332
$ws.Cells.Item(80, 11).Value = 2379099.5
$ws.Cells.Item(80, 12).Value = 34499.332
$ws.Cells.Item(80, 13).Value = -2378101.5
$ws.Cells.Item(80, 14).Value = -36495.332

$ws.Cells.Item(83, 8).Value = 1296976.4
$ws.Cells.Item(83, 9).Value = 2379099.5
$ws.Cells.Item(83, 10).Value = 34499.332
$ws.Cells.Item(83, 11).Value = 11895497.5
$ws.Cells.Item(83, 12).Value = 172496.66
$ws.Cells.Item(83, 13).Value = -11890505.5
$ws.Cells.Item(83, 14).Value = -182480.66

$ws.Cells.Item(97, 8).Value = 675.1905
$ws.Cells.Item(97, 9).Value = 541.8125
$ws.Cells.Item(97, 11).Value = 541.8125
$ws.Cells.Item(97, 13).Value = -45.8125

$ws.Cells.Item(113, 8).Value = 4349.75
$ws.Cells.Item(113, 9).Value = 3949.5
$ws.Cells.Item(113, 11).Value = 3949.5
$ws.Cells.Item(113, 13).Value = -1779.5

$ws.Cells.Item(132, 8).Value = 391582.16
$ws.Cells.Item(132, 9).Value = 127973.94
$ws.Cells.Item(132, 11).Value = 383921.82
$ws.Cells.Item(132, 13).Value = -381391.82

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 872.05884
$ws.Cells.Item(22, 10).Value = 812.5
$ws.Cells.Item(22, 12).Value = 812.5
$ws.Cells.Item(22, 14).Value = -1402.5

$ws.Cells.Item(27, 8).Value = 872.05884
$ws.Cells.Item(27, 10).Value = 812.5
$ws.Cells.Item(27, 12).Value = 812.5
$ws.Cells.Item(27, 14).Value = -1026.5

$ws.Cells.Item(46, 8).Value = 5356.8823
$ws.Cells.Item(46, 9).Value = 3866
$ws.Cells.Item(46, 11).Value = 3866
$ws.Cells.Item(46, 13).Value = -3678

$ws.Cells.Item(55, 8).Value = 320.52942
$ws.Cells.Item(55, 10).Value = 417.8
$ws.Cells.Item(55, 12).Value = 417.8
$ws.Cells.Item(55, 14).Value = -763.8

$ws.Cells.Item(122, 8).Value = 7054.8696
$ws.Cells.Item(122, 9).Value = 3388.0667
$ws.Cells.Item(122, 10).Value = 13930.125
$ws.Cells.Item(122, 11).Value = 10164.2001
$ws.Cells.Item(122, 12).Value = 41790.375
$ws.Cells.Item(122, 13).Value = -7714.2001
$ws.Cells.Item(122, 14).Value = -46690.375
